$d = $word.ActiveDocument
$r = $d.Content
$r.Collapse(0)

# Insert both new paragraphs in a single InsertXML call so the trailing
# empty-story paragraph is consumed exactly once.
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="ＭＳ 明朝" w:hAnsi="ＭＳ 明朝" w:cs="ＭＳ 明朝" w:eastAsia="ＭＳ 明朝"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:shd w:fill="auto" w:val="clear"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="ＭＳ 明朝" w:hAnsi="ＭＳ 明朝" w:cs="ＭＳ 明朝" w:eastAsia="ＭＳ 明朝"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:shd w:fill="auto" w:val="clear"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="ＭＳ 明朝" w:hAnsi="ＭＳ 明朝" w:cs="ＭＳ 明朝" w:eastAsia="ＭＳ 明朝"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:shd w:fill="auto" w:val="clear"/></w:rPr><w:t xml:space="preserve">1</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="ＭＳ 明朝" w:hAnsi="ＭＳ 明朝" w:cs="ＭＳ 明朝" w:eastAsia="ＭＳ 明朝"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:shd w:fill="auto" w:val="clear"/></w:rPr><w:t xml:space="preserve">行追加した</w:t></w:r></w:p>'
$r.InsertXML($xml) | Out-Null

# Force explicit (non-default-dropped) spacing/indent values on both new
# paragraphs to match the original document's authored XML.
$p1 = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$p1.LeftIndent = 0
$p1.RightIndent = 0
$p1.FirstLineIndent = 0
$p1.SpaceBefore = 0
$p1.SpaceAfter = 10
$p1.LineSpacing = 13.8

$p2 = $d.Paragraphs.Last
$p2.LeftIndent = 0
$p2.RightIndent = 0
$p2.FirstLineIndent = 0
$p2.SpaceBefore = 0
$p2.SpaceAfter = 10
$p2.LineSpacing = 13.8
